$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A19").Value = "7 Inch Combo Pack"
$ws.Range("C19").Value = 350
$ws.Range("C19").NumberFormat = $ws.Range("C18").NumberFormat
$ws.Range("D19").Value = "7 Inch Combo Pack.png"

$ws.Range("J11").Select()
